$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for the "Poroto verde" series at
# "Feria Lagunitas de Puerto Montt". It belongs right above the existing
# row 38 (chronologically it is the most recent entry), so insert a new
# row there and push every following record down by one - this is exactly
# what happened to rows 38..91 (now 39..92) in the target workbook.
$ws.Rows("38:38").Insert()

$ws.Range("A38").Value = 4
$ws.Range("B38").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C38").Value = "Los Lagos"
$ws.Range("D38").Value = 44771
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = 100112031
$ws.Range("G38").Value = "Poroto verde"
$ws.Range("H38").Value = "Magnum"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 35
$ws.Range("K38").Value = 40000
$ws.Range("L38").Value = 40000
$ws.Range("M38").Value = 40000
$ws.Range("N38").Value = "$/malla 25 kilos"
$ws.Range("O38").Value = "Perú"
$ws.Range("P38").Value = 1600
$ws.Range("Q38").Value = 25
$ws.Range("R38").Value = "Hortaliza"
